# Applies the "Deploying to gh-pages" content refresh to the
# StructureDefinition-based-on-value workbook:
#   - Metadata sheet: bump Version / Date, fill in Publisher, replace the
#     duplicated "Contact" rows with a single "Jurisdiction" row, which
#     shifts every row below it up by one.
#   - Elements sheet: the root Extension row's Short/Definition text is
#     replaced with the profile's own Title/Description.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bumped to the new publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a second "Publisher" detail row that said
# "Contact" / "No display for ContactDetail". It becomes the
# "Jurisdiction" / "United States of America" row instead.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was an exact duplicate of the old row 10 ("Contact" /
# "No display for ContactDetail"); delete it outright so every
# subsequent row (Description, Purpose, Copyright, ...) shifts up by one
# and the sheet's used range shrinks from A1:B21 to A1:B20.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# The root Extension element's Short/Definition columns now surface the
# profile's own Title/Description instead of the generic placeholders.
$elements.Range("K2").Value = "Based On Value"
$elements.Range("L2").Value = "Value the insight result used as input"
